# CDL_limit.xlsx - "Add files via upload" edit
# Updates several CDL sensor limit rows on Sheet1: new min/max bounds,
# a couple of units switched from kPa to psi, and the ground-speed row's
# name normalized from "cdlgroundspeed" to "CDLGroundSpeed". Also moves
# the active selection to D5 (scrolled back to the top of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 24 (CDLGroundSpeed / was "cdlgroundspeed"): fix the capitalization of
# the metric's name. Do this before the other edits below so the new shared
# strings ("CDLGroundSpeed" then "psi") land in the same order Excel wrote
# them in the reference edit.
$ws.Range("A24").Value = "CDLGroundSpeed"

# Row 5 (CDLAtmosphericPressure): widen the range and report it in psi.
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 4750
$ws.Range("D5").Value = "psi"

# Row 9 (CDLBrakeFilterBypassStatus): raise max_value.
$ws.Range("C9").Value = 100

# Row 14 (CDLEngineOilPressureAbs): raise max_value.
$ws.Range("C14").Value = 600

# Row 20 (CDLFuelConsumptionRateLPH): raise max_value.
$ws.Range("C20").Value = 200

# Row 50 (CDLSteeringPumpOilPressure): raise max_value.
$ws.Range("C50").Value = 10000

# Row 57 (CDLTransmissionOilPressureAbs): widen the range and report it in psi.
$ws.Range("B57").Value = 0
$ws.Range("C57").Value = 4000
$ws.Range("D57").Value = "psi"

# Row 63 (CDLTurbocharger3CompressorInletPressure): lower min_value.
$ws.Range("B63").Value = 10

# Move the selection to D5 and scroll the view back to the top (the
# original file had scrolled down to row 35 and left B65 selected).
$ws.Activate()
$null = $ws.Range("D5").Select()
